# Delete the sample row that had no fastq file in the 06.17.19 library.
# This corresponds to the row with s2cDNASampleNumber = 12 (worksheet row 13),
# which is removed via a native row delete (cells below shift up one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$null = $ws.Range("A13").EntireRow.Delete()

# Move the selection to match the post-edit cursor position.
$null = $ws.Range("E10").Select()
